$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = "67.741.70"
$ws.Cells.Item(2, 5).Value = "  +3.14%  "

# Row 3
$ws.Cells.Item(3, 4).Value = "3.320.16"
$ws.Cells.Item(3, 5).Value = "  +0.67%  "

# Row 4
$ws.Cells.Item(4, 4).Value = "'1.00"
$ws.Cells.Item(4, 5).Value = "  -0.08%  "

# Row 5
$ws.Cells.Item(5, 4).Value = "'585.31"
$ws.Cells.Item(5, 5).Value = "  +5.36%  "

# Row 6
$ws.Cells.Item(6, 4).Value = "'181.89"
$ws.Cells.Item(6, 5).Value = "  -1.04%  "

# Row 7
$ws.Cells.Item(7, 4).Value = "'1.00"
$ws.Cells.Item(7, 5).Value = "  -0.11%  "

# Row 8
$ws.Cells.Item(8, 5).Value = "  +3.23%  "

# Row 9
$ws.Cells.Item(9, 4).Value = "3.311.15"
$ws.Cells.Item(9, 5).Value = "  +0.63%  "

# Row 10
$ws.Cells.Item(10, 5).Value = "  +3.94%  "

# Row 11
$ws.Cells.Item(11, 4).Value = "'0.579"
$ws.Cells.Item(11, 5).Value = "  +0.97%  "

# Row 12
$ws.Cells.Item(12, 4).Value = "'46.38"
$ws.Cells.Item(12, 5).Value = "  +2.27%  "

# Row 13
$ws.Cells.Item(13, 4).Value = "'0.0000276"
$ws.Cells.Item(13, 5).Value = "  +6.75%  "

# Row 14
$ws.Cells.Item(14, 4).Value = "'638.23"
$ws.Cells.Item(14, 5).Value = "  +11.74%  "

# Row 15
$ws.Cells.Item(15, 4).Value = "3.858.20"
$ws.Cells.Item(15, 5).Value = "  +0.76%  "

# Row 16
$ws.Cells.Item(16, 4).Value = "'8.44"
$ws.Cells.Item(16, 5).Value = "  +1.09%  "

# Row 17
$ws.Cells.Item(17, 4).Value = "67.892.85"
$ws.Cells.Item(17, 5).Value = "  +3.51%  "

# Row 18
$ws.Cells.Item(18, 5).Value = "  +1.52%  "

# Row 19
$ws.Cells.Item(19, 4).Value = "3.319.47"
$ws.Cells.Item(19, 5).Value = "  +0.76%  "

# Row 20
$ws.Cells.Item(20, 4).Value = "'17.67"
$ws.Cells.Item(20, 5).Value = "  +0.79%  "

# Row 21
$ws.Cells.Item(21, 4).Value = "'10.95"
$ws.Cells.Item(21, 5).Value = "  +2.18%  "

# Row 22
$ws.Cells.Item(22, 4).Value = "'0.898"
$ws.Cells.Item(22, 5).Value = "  +1.68%  "

# Row 23
$ws.Cells.Item(23, 4).Value = "'17.72"
$ws.Cells.Item(23, 5).Value = "  +0.15%  "

# Row 24
$ws.Cells.Item(24, 4).Value = "'5.02"
$ws.Cells.Item(24, 5).Value = "  +1.09%  "

# Row 25
$ws.Cells.Item(25, 4).Value = "'97.87"
$ws.Cells.Item(25, 5).Value = "  -0.30%  "

# Row 26
$ws.Cells.Item(26, 4).Value = "'3.99"
$ws.Cells.Item(26, 5).Value = "  +1.97%  "

# Row 27
$ws.Cells.Item(27, 4).Value = "'2.80"
$ws.Cells.Item(27, 5).Value = "  +5.28%  "

# Row 28
$ws.Cells.Item(28, 4).Value = "'9.61"
$ws.Cells.Item(28, 5).Value = "  +4.03%  "

# Row 29
$ws.Cells.Item(29, 4).Value = "'32.80"
$ws.Cells.Item(29, 5).Value = "  +8.38%  "

# Row 30
$ws.Cells.Item(30, 4).Value = "'8.57"
$ws.Cells.Item(30, 5).Value = "  +2.08%  "

# Row 31
$ws.Cells.Item(31, 4).Value = "'6.66"
$ws.Cells.Item(31, 5).Value = "  +0.99%  "

# Row 32
$ws.Cells.Item(32, 4).Value = "'605.15"
$ws.Cells.Item(32, 5).Value = "  +9.07%  "

# Row 33
$ws.Cells.Item(33, 4).Value = "3.929.91"
$ws.Cells.Item(33, 5).Value = "  +4.75%  "

# Row 34
$ws.Cells.Item(34, 2).Value = "Cosmos"
$ws.Cells.Item(34, 3).Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Cells.Item(34, 4).Value = "'10.95"
$ws.Cells.Item(34, 5).Value = "  +1.94%  "

# Row 35
$ws.Cells.Item(35, 2).Value = "dogwifhat"
$ws.Cells.Item(35, 3).Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Cells.Item(35, 4).Value = "'3.60"
$ws.Cells.Item(35, 5).Value = "  -1.05%  "

# Row 36
$ws.Cells.Item(36, 5).Value = "  +2.57%  "

# Row 37
$ws.Cells.Item(37, 4).Value = "'0.999"
$ws.Cells.Item(37, 5).Value = "  +0.01%  "

# Row 38
$ws.Cells.Item(38, 4).Value = "'55.74"
$ws.Cells.Item(38, 5).Value = "  +0.54%  "

# Row 39
$ws.Cells.Item(39, 4).Value = "'3.27"
$ws.Cells.Item(39, 5).Value = "  +5.70%  "

# Row 40
$ws.Cells.Item(40, 2).Value = "Kaspa"
$ws.Cells.Item(40, 3).Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Cells.Item(40, 4).Value = "'0.128"
$ws.Cells.Item(40, 5).Value = "  +2.40%  "

# Row 41
$ws.Cells.Item(41, 2).Value = "Fetch.AI"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Cells.Item(41, 4).Value = "'2.70"
$ws.Cells.Item(41, 5).Value = "  +6.61%  "

# Row 42
$ws.Cells.Item(42, 4).Value = "'32.96"

# Row 43
$ws.Cells.Item(43, 4).Value = "0.0₃0688"
$ws.Cells.Item(43, 5).Value = "  +2.98%  "

# Row 44
$ws.Cells.Item(44, 5).Value = "  +1.00%  "

# Row 45
$ws.Cells.Item(45, 4).Value = "'0.337"
$ws.Cells.Item(45, 5).Value = "  +2.80%  "

# Row 46
$ws.Cells.Item(46, 5).Value = "  +2.77%  "

# Row 47
$ws.Cells.Item(47, 4).Value = "'0.128"
$ws.Cells.Item(47, 5).Value = "  +2.61%  "

# Row 48
$ws.Cells.Item(48, 4).Value = "'1.01"
$ws.Cells.Item(48, 5).Value = "  +0.68%  "

# Row 49
$ws.Cells.Item(49, 4).Value = "'2.54"
$ws.Cells.Item(49, 5).Value = "  +2.51%  "

# Row 50
$ws.Cells.Item(50, 5).Value = "  +9.31%  "

# Row 51
$ws.Cells.Item(51, 4).Value = "'130.98"
$ws.Cells.Item(51, 5).Value = "  +4.69%  "
